# feat: add 2022-Q4 data
#
# The workbook tracked quarterly fund-holding data. This adds a new
# "2022-Q4" worksheet (placed right after "总计"/Totals, before the
# existing "2022-Q3" sheet) and updates the "总计" summary sheet so its
# most-recent row reflects Q4 while the previous Q3 totals move down to
# their own row.

$wb = $excel.ActiveWorkbook
$totals = $wb.Worksheets.Item(1)

# --- Insert the new "2022-Q4" sheet right after "总计" ---------------------
# Copying "总计" gives the new sheet the same header/number styling and
# page margins that the rest of the "2022-Q4" tab should have, then we
# overwrite its contents with the quarterly fund-detail data.
$totals.Copy($null, $totals)
$q4 = $wb.Worksheets.Item(2)
$q4.Name = "2022-Q4"

# The totals sheet only had styled cells in B1:D1 — stretch that same
# header style across the rest of the header row (E1:H1).
$q4.Range("D1").Copy($q4.Range("E1:H1"))

$q4.Range("B1").Value = "基金代码"
$q4.Range("C1").Value = "基金名称"
$q4.Range("D1").Value = "基金规模"
$q4.Range("E1").Value = "股票总仓位"
$q4.Range("F1").Value = "仓位占比"
$q4.Range("G1").Value = "持有市值(亿元)"
$q4.Range("H1").Value = "仓位排名"

$q4.Range("A2").Value = 0

# These look like numbers but must stay text, matching the source data —
# format as text first, then clear the formatting back off afterwards so
# no stray number-format styling is left applied to the cells.
$q4.Range("B2:G2").NumberFormat = "@"
$q4.Range("B2").Value = "590003"
$q4.Range("C2").Value = "中邮核心优势灵活配置混合"
$q4.Range("D2").Value = "19.01"
$q4.Range("E2").Value = "79.84"
$q4.Range("F2").Value = "5.19"
$q4.Range("G2").Value = "0.9866"
$q4.Range("B2:G2").ClearFormats()

$q4.Range("H2").Value = 6

# --- Update the "总计" (totals) sheet --------------------------------------
# Row 2 keeps its style/format but now reports the new 2022-Q4 totals;
# the old 2022-Q3 totals move down to a new row 3 (cloning A2's style
# for A3 first, then overwriting the value).
$totals.Range("A2").Copy($totals.Range("A3"))
$totals.Range("A3").Value = 1
$totals.Range("B3").Value = "2022-Q3"
$totals.Range("C3").Value = 1
$totals.Range("D3").Value = 0.96

$totals.Range("B2").Value = "2022-Q4"
$totals.Range("D2").Value = 0.99

# Keep "2022-Q3" as the selected tab, matching the original workbook state.
$wb.Worksheets.Item("2022-Q3").Activate()
